# Claimed another two sensors
# Adds developer/finished/mode/interface info for the MagneticSensor (row 44)
# and RCXLightSensor (row 61) rows, and fills in the "Fits in framework"
# column (C) for the remaining RCX rows (62-68) on the "Blad1" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Blad1")

# Row 44 - MagneticSensor: HiTechnicMagneticSensor
$ws.Cells.Item(44, 6).Value = "Magnetic"                   # F44 - Mode
$ws.Cells.Item(44, 2).Value = "HiTechnicMagneticSensor"   # B44 - new name
$ws.Cells.Item(44, 4).Value = "Lawrie"                     # D44 - Developer
$ws.Cells.Item(44, 5).Value = "N"                          # E44 - finished
$ws.Cells.Item(44, 7).Value = "SampleProvider"             # G44 - interface

# Row 61 - RCXLightSensor (new name same as old name)
$ws.Cells.Item(61, 2).Value = "RCXLightSensor"              # B61 - new name
$ws.Cells.Item(61, 4).Value = "Lawrie"                     # D61 - Developer
$ws.Cells.Item(61, 5).Value = "N"                          # E61 - finished
$ws.Cells.Item(61, 6).Value = "Light"                      # F61 - Mode
$ws.Cells.Item(61, 7).Value = "SampleProvider"             # G61 - interface

# Rows 62-66, 68 - "Fits in framework" = N
$ws.Cells.Item(62, 3).Value = "N"
$ws.Cells.Item(63, 3).Value = "N"
$ws.Cells.Item(64, 3).Value = "N"
$ws.Cells.Item(65, 3).Value = "N"
$ws.Cells.Item(66, 3).Value = "N"
# Row 67 - "Fits in framework" = ?
$ws.Cells.Item(67, 3).Value = "?"
$ws.Cells.Item(68, 3).Value = "N"

$ws.Range("C68").Select() | Out-Null
